# Fraction of Forests Owned by Entity.xlsx — update to v2.0.0 content
# (adds "Biomass" ownership notes on About, splits FoFObE categories into
#  the full EPS entity list, and applies header/label formatting.)

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")
$wsFoF   = $wb.Worksheets.Item("FoFObE")

# ---------------------------------------------------------------------
# About sheet: append a new "Biomass" notes block below the existing
# "Notes" section (rows 15-22; row 14 left blank like row 8).
# ---------------------------------------------------------------------
$wsAbout.Range("A15").Value = "Biomass"
$wsAbout.Range("A15").Font.Bold = $true

$wsAbout.Range("A16").Value = "Though some industry-owned timber land might be used to"
$wsAbout.Range("A17").Value = "produce biomass, but not all biomass is wood, and"
$wsAbout.Range("A18").Value = "much of the wood biomass is waste (bark, sawdust, chips, scrap,"
$wsAbout.Range("A19").Value = "and paper mill residues), which is not the main product"
$wsAbout.Range("A20").Value = "of the timber lands.  Accordingly, we assign timber land"
$wsAbout.Range("A21").Value = "ownership to ""nonenergy industries"" rather than assigning"
$wsAbout.Range("A22").Value = "a share to ""biomass and biofuel suppliers."""

# ---------------------------------------------------------------------
# Data sheet: the "U.S." header label (B2) becomes bold + right aligned.
# ---------------------------------------------------------------------
$wsData.Range("B2").Font.Bold = $true
$wsData.Range("B2").HorizontalAlignment = -4152   # xlRight
$wsData.PageSetup.Orientation = 1                  # xlPortrait

# ---------------------------------------------------------------------
# FoFObE sheet: relabel the header, rename two categories and add the
# six new zero-valued ownership categories used by later EPS versions.
# ---------------------------------------------------------------------
$wsFoF.Range("B1").Value = "Fraction of Forest Owned (dimensionless)"
$wsFoF.Range("B1").HorizontalAlignment = -4152     # xlRight
$wsFoF.Range("B1").WrapText = $true
$wsFoF.Rows.Item(1).RowHeight = 28.5

$wsFoF.Range("A3").Value = "nonenergy industries"
$wsFoF.Range("A4").Value = "labor and consumers"

$wsFoF.Range("A5").Value = "foreign entities"
$wsFoF.Range("B5").Value = 0

$wsFoF.Range("A6").Value = "electricity suppliers"
$wsFoF.Range("B6").Value = 0

$wsFoF.Range("A7").Value = "coal suppliers"
$wsFoF.Range("B7").Value = 0

$wsFoF.Range("A8").Value = "natural gas and petroleum suppliers"
$wsFoF.Range("B8").Value = 0

$wsFoF.Range("A9").Value = "biomass and biofuel suppliers"
$wsFoF.Range("B9").Value = 0

$wsFoF.Range("A10").Value = "other energy suppliers"
$wsFoF.Range("B10").Value = 0

$wsFoF.Columns.Item(1).ColumnWidth = 34

# Keep "About" as the visible/active sheet & selection, but remember the
# FoFObE sheet's own B1 selection (matches the saved file's per-sheet view).
$wsFoF.Range("B1").Select()
$wsAbout.Activate()
$wsAbout.Range("A1").Select()
